# Auto-generated Excel COM-interop script to apply the Hyperion_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5723.75
$ws.Range("J51").Value = 5165
$ws.Range("L51").Value = 5165
$ws.Range("N51").Value = -6133
$ws.Range("H138").Value = 3243.2407
$ws.Range("J138").Value = 3869.3057
$ws.Range("L138").Value = 11607.9171
$ws.Range("N138").Value = -21887.9171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1969.6
$ws.Range("J2").Value = 1954.3334
$ws.Range("L2").Value = 1954.3334
$ws.Range("N2").Value = -2180.3334
$ws.Range("H32").Value = 4151.2256
$ws.Range("I32").Value = 2060.35
$ws.Range("K32").Value = 2060.35
$ws.Range("M32").Value = -1773.35
$ws.Range("H45").Value = 53130.9
$ws.Range("I45").Value = 92415.09
$ws.Range("K45").Value = 92415.09
$ws.Range("M45").Value = -92038.09
$ws.Range("H63").Value = 2024.4166
$ws.Range("I63").Value = 1853.909
$ws.Range("J63").Value = 3900
$ws.Range("K63").Value = 1853.909
$ws.Range("L63").Value = 3900
$ws.Range("M63").Value = -1167.909
$ws.Range("N63").Value = -5272
$ws.Range("H66").Value = 2024.4166
$ws.Range("I66").Value = 1853.909
$ws.Range("J66").Value = 3900
$ws.Range("K66").Value = 9269.545
$ws.Range("L66").Value = 19500
$ws.Range("M66").Value = -5837.545
$ws.Range("N66").Value = -26364
$ws.Range("H74").Value = 41828.98
$ws.Range("I74").Value = 5879.324
$ws.Range("J74").Value = 152673.75
$ws.Range("K74").Value = 5879.324
$ws.Range("L74").Value = 152673.75
$ws.Range("M74").Value = -5005.324
$ws.Range("N74").Value = -154421.75
$ws.Range("H77").Value = 41828.98
$ws.Range("I77").Value = 5879.324
$ws.Range("J77").Value = 152673.75
$ws.Range("K77").Value = 29396.62
$ws.Range("L77").Value = 763368.75
$ws.Range("M77").Value = -25028.62
$ws.Range("N77").Value = -772104.75
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H110").Value = 4765.5
$ws.Range("I110").Value = 3660.6428
$ws.Range("K110").Value = 3660.6428
$ws.Range("M110").Value = -1615.6428
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H116").Value = 1969.6
$ws.Range("J116").Value = 1954.3334
$ws.Range("L116").Value = 1954.3334
$ws.Range("N116").Value = -6542.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1969.6
$ws.Range("J3").Value = 1954.3334
$ws.Range("L3").Value = 1954.3334
$ws.Range("N3").Value = -2182.3334
$ws.Range("H80").Value = 292.1579
$ws.Range("J80").Value = 342.25
$ws.Range("L80").Value = 342.25
$ws.Range("N80").Value = -2338.25
$ws.Range("H83").Value = 292.1579
$ws.Range("J83").Value = 342.25
$ws.Range("L83").Value = 1711.25
$ws.Range("N83").Value = -11695.25
$ws.Range("H94").Value = 8463.333
$ws.Range("I94").Value = 1695.3334
$ws.Range("J94").Value = 15231.333
$ws.Range("K94").Value = 1695.3334
$ws.Range("L94").Value = 15231.333
$ws.Range("M94").Value = -1244.3334
$ws.Range("N94").Value = -16133.333
$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -53134
$ws.Range("H134").Value = 4647.7856
$ws.Range("I134").Value = 2441.5
$ws.Range("K134").Value = 7324.5
$ws.Range("M134").Value = -4789.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 38098
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 1899.5358
$ws.Range("I58").Value = 1869.9231
$ws.Range("J58").Value = 1925.2
$ws.Range("K58").Value = 1869.9231
$ws.Range("L58").Value = 1925.2
$ws.Range("M58").Value = -1666.9231
$ws.Range("N58").Value = -2331.2
$ws.Range("H61").Value = 38098
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H107").Value = 2579.8
$ws.Range("I107").Value = 2586.182
$ws.Range("K107").Value = 2586.182
$ws.Range("M107").Value = -666.1819999999998
$ws.Range("H132").Value = 71576.46
$ws.Range("I132").Value = 3088.25
$ws.Range("J132").Value = 181157.6
$ws.Range("K132").Value = 9264.75
$ws.Range("L132").Value = 543472.8
$ws.Range("M132").Value = -6734.75
$ws.Range("N132").Value = -548532.8
$ws.Range("H136").Value = 1899.5358
$ws.Range("I136").Value = 1869.9231
$ws.Range("J136").Value = 1925.2
$ws.Range("K136").Value = 5609.7693
$ws.Range("L136").Value = 5775.6
$ws.Range("M136").Value = -3059.7693
$ws.Range("N136").Value = -10875.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1003.67645
$ws.Range("J5").Value = 1424.7858
$ws.Range("L5").Value = 4274.357400000001
$ws.Range("N5").Value = -4498.357400000001
$ws.Range("H26").Value = 306
$ws.Range("I26").Value = 374.66666
$ws.Range("K26").Value = 1123.99998
$ws.Range("M26").Value = -835.99998
$ws.Range("H36").Value = 200
$ws.Range("J36").Value = 200
$ws.Range("L36").Value = 600
$ws.Range("N36").Value = -938
$ws.Range("H39").Value = 3042.8572
$ws.Range("I39").Value = 950
$ws.Range("J39").Value = 3880
$ws.Range("K39").Value = 2850
$ws.Range("L39").Value = 11640
$ws.Range("M39").Value = -2556
$ws.Range("N39").Value = -12228
$ws.Range("H56").Value = 20838824
$ws.Range("I56").Value = 20838824
$ws.Range("K56").Value = 20838824
$ws.Range("M56").Value = -20838294
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H129").Value = 2501310.5
$ws.Range("J129").Value = 3033
$ws.Range("L129").Value = 9099
$ws.Range("N129").Value = -19099
$ws.Range("H131").Value = 19844572
$ws.Range("J131").Value = 20837482
$ws.Range("L131").Value = 62512446
$ws.Range("N131").Value = -62522526
$ws.Range("H132").Value = 1566.0322
$ws.Range("I132").Value = 1080.6154
$ws.Range("J132").Value = 1916.6111
$ws.Range("K132").Value = 9725.5386
$ws.Range("L132").Value = 17249.4999
$ws.Range("M132").Value = -7195.5386
$ws.Range("N132").Value = -22309.4999
$ws.Range("H135").Value = 1003.67645
$ws.Range("J135").Value = 1424.7858
$ws.Range("L135").Value = 12823.0722
$ws.Range("N135").Value = -17893.0722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 39999.4
$ws.Range("I124").Value = 39999
$ws.Range("J124").Value = 39999.5
$ws.Range("K124").Value = 39999
$ws.Range("L124").Value = 39999.5
$ws.Range("M124").Value = -35089
$ws.Range("N124").Value = -49819.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 8000
$ws.Range("J28").Value = 8000
$ws.Range("L28").Value = 8000
$ws.Range("N28").Value = -8696
$ws.Range("H31").Value = 7000
$ws.Range("J31").Value = 7000
$ws.Range("L31").Value = 7000
$ws.Range("N31").Value = -7696
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H96").Value = 1568.8
$ws.Range("I96").Value = 1495
$ws.Range("J96").Value = 1587.25
$ws.Range("K96").Value = 1495
$ws.Range("L96").Value = 1587.25
$ws.Range("M96").Value = -122
$ws.Range("N96").Value = -4333.25
$ws.Range("H107").Value = 3912.4285
$ws.Range("I107").Value = 4645.909
$ws.Range("K107").Value = 13937.727
$ws.Range("M107").Value = -12017.727
